$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The text value in D2 was misspelled as "NofInstalments" -> corrected to "Noofinstalments"
$ws.Range("D2").Value = "Noofinstalments"

# Column D widened slightly to fit the longer corrected text
$ws.Columns("D").ColumnWidth = 13.6667

# Move / update the active selection to Q2 (matches the author's last selection when saving)
$ws.Range("Q2").Select()
